# "addition col to reflect KO Bins"
# Adds a new "stateOfMatter" column (F) for the existing rows (gold/silver/
# platinum/argon/boron), and adds a brand-new row (7) describing mercury,
# including its stateOfMatter value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F: stateOfMatter -------------------------------------------------
$ws.Range("F1").Value = "stateOfMatter"
$ws.Range("F2").Value = "solid"
$ws.Range("F3").Value = "solid"
$ws.Range("F4").Value = "solid"
$ws.Range("F5").Value = "gas"
$ws.Range("F6").Value = "gas"

# --- New row 7: mercury ------------------------------------------------------
$ws.Range("A7").Value = "mercury"
$ws.Range("B7").Value = "Mercury is a chemical element with symbol Hg and atomic number 80. Classified as a transition metal, Mercury is a liquid at room temperature."
$ws.Range("C7").Value = "[Xe]6s24f145d10"
$ws.Range("D7").Value = "Hg"
$ws.Range("E7").Value = "The name derives from the Roman god Mercury, the nimble messenger of the gods, because the ancients used that name for the element known from prehistoric times. The symbol Hg derives from the Greek hydrargyrum for `"liquid silver`" or `"quick silver`"."
$ws.Range("F7").Value = "liquid"

# Wrap the long etymology text in the new mercury row, matching the wrapped
# etymology cells already used for the other elements.
$ws.Range("E7").WrapText = $true

# --- Row heights / column widths --------------------------------------------
$ws.Rows("1:7").RowHeight = 98.5
$ws.Columns("B").ColumnWidth = 43.26
$ws.Columns("E").ColumnWidth = 52.92

# --- Selection / view ---------------------------------------------------------
$ws.Range("G8").Select()
